{"js": "// \"removed year on all title slides\"\n// The title-slide front matter is: Title paragraph, Author paragraph,\n// then a Date paragraph (e.g. \"2020-11-25\"). Remove every paragraph\n// styled \"Date\" (the year/date line under the author) from the document.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.style === \"Date\") {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"removed year on all title slides\"\n# The title-slide front matter is: Title paragraph, Author paragraph,\n# then a Date paragraph (e.g. \"2020-11-25\"). Remove every paragraph\n# styled \"Date\" (the year/date line under the author) from the document.\n\n$d = $word.ActiveDocument\n\n# Walk backwards so deleting a paragraph doesn't shift the indices of\n# paragraphs we still need to visit.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Date\") {\n        $p.Range.Delete()\n    }\n}\n"}
